$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IMPORT_DATA")
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
